# Update "想去人数" (want-to-go count) figures in column F across sheets,
# reflecting the refreshed data output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 101
$ws1.Range("F4").Value = 614
$ws1.Range("F6").Value = 9249
$ws1.Range("F7").Value = 833
$ws1.Range("F9").Value = 1183
$ws1.Range("F10").Value = 1079
$ws1.Range("F11").Value = 139
$ws1.Range("F12").Value = 70
$ws1.Range("F14").Value = 255
$ws1.Range("F15").Value = 377
$ws1.Range("F18").Value = 1214

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 8

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 101
$ws4.Range("F5").Value = 8
$ws4.Range("F6").Value = 614
$ws4.Range("F8").Value = 9249
$ws4.Range("F9").Value = 833
$ws4.Range("F11").Value = 1183
$ws4.Range("F12").Value = 1079
$ws4.Range("F13").Value = 139
$ws4.Range("F14").Value = 70
$ws4.Range("F16").Value = 255
$ws4.Range("F17").Value = 377
$ws4.Range("F20").Value = 1214
